$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 79, shifting existing rows 79:84 down to 80:85
$ws.Rows("79:79").Insert()

# Fill the new row 79 with the new data record (same categorical fields as the
# other "Vega Modelo de Temuco - Tuna" rows, new date/quality/volume/price data)
$ws.Range("A79").Value = 10
$ws.Range("B79").Value = "Vega Modelo de Temuco"
$ws.Range("C79").Value = "La Araucanía"
$ws.Range("D79").Value = 45013
$ws.Range("E79").Value = 9
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100107
$ws.Range("H79").Value = "Otros"
$ws.Range("I79").Value = 100107011
$ws.Range("J79").Value = "Tuna"
$ws.Range("K79").Value = "Sin especificar"
$ws.Range("L79").Value = "Primera"
$ws.Range("M79").Value = 45
$ws.Range("N79").Value = 18000
$ws.Range("O79").Value = 18000
$ws.Range("P79").Value = 18000
$ws.Range("Q79").Value = "$/caja 16 kilos"
$ws.Range("R79").Value = "Provincia de Los Andes"
$ws.Range("S79").Value = 1125
$ws.Range("T79").Value = 16
